$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(2)

# Before: "12" + "th"(sup) + " or 26" + "th"(sup) + " of August"
# After:  "12" + "th"(sup) + " of " + "August"
$run3 = $para.Runs(3)
$run3.Text = " of "

# Remove the now-redundant superscript "th" run that followed "26"
$run4 = $para.Runs(4)
$run4.Text = ""

# The trailing run's text changes from " of August" to "August"
$run5 = $para.Runs(4)
$run5.Text = "August"

# Text removal can trigger the textbox's auto-fit height recalculation;
# restore the original shape height since the visible line count is unchanged.
$shp.Height = 1195199 / 12700
